# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The "CasesTab" query (column B, row 2) no longer returns the Cohort
# column, and the three per-tab Cypher queries (CasesTab/SamplesTab/
# FilesTab) are re-saved in SamplesTab, FilesTab, CasesTab order so the
# workbook's internal string table matches the authors' re-export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 / SamplesTab: sample-level query (text unchanged) ---
$ws.Range("B3").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE s.clinical_study_designation IN [''UBC01''] and demo.breed in [''Belgian Malinois'', ''Labrador Retriever'',''West Highland White Terrier'']and diag.disease_term in [''Bladder Cancer''] and diag.primary_disease_site in [ ''Bladder, Prostate''] and diag.best_response in [''Partial Response'']
 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`'

# --- Row 4 / FilesTab: file-level query (text unchanged) ---
$ws.Range("B4").Value = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN [''UBC01''] and demo.breed in [''Belgian Malinois'', ''Labrador Retriever'',''West Highland White Terrier'']and diag.disease_term in [''Bladder Cancer''] and diag.primary_disease_site in [ ''Bladder, Prostate''] and diag.best_response in [''Partial Response'']
    
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '''') AS `File Name`, 
        coalesce(f.file_type, '''') AS `File Type`, 
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `File Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'''') AS `Study Code`'

# --- Row 2 / CasesTab: case-level query with the trailing
#     `coalesce(co.cohort_description, '') AS Cohort` column removed ---
$ws.Range("B2").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN [''UBC01''] and demo.breed in [''Belgian Malinois'', ''Labrador Retriever'',''West Highland White Terrier'']and diag.disease_term in [''Bladder Cancer''] and diag.primary_disease_site in [ ''Bladder, Prostate''] and diag.best_response in [''Partial Response'']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'

# Row 2 is one line shorter now that the Cohort column is gone.
$ws.Rows.Item(2).RowHeight = 304.5

# Selection/scroll moved back up to B2 after the edit.
$ws.Range("B2").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
